$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (D, J, K, L, M, N, O, P) resulting from reshuffled dataset
$rows = @{
    2  = @{ D = 44169; J = 3000; K = 1000; L = 1000; M = 1000; N = '$/kilo';  O = 'Provincia de Linares'; P = 1000 }
    3  = @{ D = 44161; J = 3000; K = 1000; L = 1000; M = 1000; N = '$/kilo';  O = 'Provincia de Linares'; P = 1000 }
    4  = @{ D = 44167; J = 2000; K = 1000; L = 1000; M = 1000; N = '$/kilo';  O = 'Región del Maule';     P = 1000 }
    6  = @{ D = 44172; J = 2000; K = 1000; L = 1000; M = 1000; N = '$/kilo';  O = 'Región del Maule';     P = 1000 }
    7  = @{ D = 44166; J = 1500; K = 1000; L = 1000; M = 1000; N = '$/kilo';  O = 'Provincia de Linares'; P = 1000 }
    9  = @{ D = 44160; J = 2000; K = 800;  L = 800;  M = 800;  N = '$/kilo';  O = 'Provincia de Linares'; P = 800  }
    10 = @{ D = 44175; J = 800;  K = 1000; L = 1100; M = 1050; N = '$/kilo';  O = 'Provincia de Linares'; P = 1050 }
    11 = @{ D = 44162; J = 4000; K = 1000; L = 1000; M = 1000; N = '$/atado'; O = 'Provincia de Linares'; P = 1000 }
    12 = @{ D = 44176; J = 2000; K = 900;  L = 900;  M = 900;  N = '$/kilo';  O = 'Provincia de Linares'; P = 900  }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
}
